$wb = $excel.ActiveWorkbook

# --- Rename the "#TParties" sheet to "#TTParties" (SIAMv4 rename: TParty -> TTParty) ---
$wsParties = $wb.Worksheets.Item("#TParties")
$wsParties.Name = "#TTParties"

# --- Update the renamed concept strings on the #TTParties sheet ---
$wsParties.Range("A1").Value = "[TTParties]"
$wsParties.Range("B1").Value = "ttIsaTTParty"

$wsParties.Range("F1").Value = "ttPartyReqdPartyRef"
$wsParties.Range("G1").Value = "ttPartyReqdPartyRef"
$wsParties.Range("H1").Value = "ttPartyReqdPartyRef"
$wsParties.Range("I1").Value = "ttPartyReqdPartyRef"

$wsParties.Range("F2").Value = "PartyRef"
$wsParties.Range("G2").Value = "PartyRef"
$wsParties.Range("H2").Value = "PartyRef"
$wsParties.Range("I2").Value = "PartyRef"

# column B goes from a wide hidden helper column to a thin visible spacer column
$wsParties.Columns.Item(2).Hidden = $false
$wsParties.Columns.Item(2).ColumnWidth = 1.21875

# page setup (portrait, paper size 9 = A4) picked up for this sheet
$wsParties.PageSetup.PaperSize = 9
$wsParties.PageSetup.Orientation = 1

# leftover cursor position after editing this sheet
$wsParties.Range("D11").Select()

# --- #Concerns sheet: objvSHRoleName -> objSHRoleName ---
$wsConcerns = $wb.Worksheets.Item("#Concerns")
$wsConcerns.Range("E1").Value = "objSHRoleName"
$wsConcerns.Range("E1").Select()

# --- #Organizations sheet: orgAbbrName/OrgAbbrName -> orgRef/OrgRef ---
$wsOrgs = $wb.Worksheets.Item("#Organizations")
$wsOrgs.Range("B1").Value = "orgRef"
$wsOrgs.Range("B2").Value = "OrgRef"
$wsOrgs.Range("B3").Select()

# --- #TTexts sheet: leftover cursor position ---
$wsTexts = $wb.Worksheets.Item("#TTexts")
$wsTexts.Range("E25").Select()

# restore the #TTParties tab as the active tab/sheet, matching activeTab="1"
$wsParties.Activate()
$wsParties.Range("D11").Select()
